# Refresh the South Korea Insurance (Prop/Cas.) capital-structure database:
# re-pull historical growth (D/E) columns, refresh every metric column for
# the already-listed companies, and re-sort rows 2-9 by company (several
# companies moved to a different row, which shows up as company_name swaps
# plus the accompanying metric refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = [double]"0.0379"  # D2
$ws.Cells.Item(2, 5).Value = [double]"0.04065000000000001"  # E2
$ws.Cells.Item(2, 6).Value = [double]"0.0318"  # F2
$ws.Cells.Item(2, 7).Value = [double]"0.02257217109468154"  # G2
$ws.Cells.Item(2, 8).Value = [double]"0.02257217109468154"  # H2
$ws.Cells.Item(2, 9).Value = [double]"0.03312017602667904"  # I2
$ws.Cells.Item(2, 10).Value = [double]"0.02805589293286049"  # J2
$ws.Cells.Item(2, 11).Value = [double]"1713.29"  # K2
$ws.Cells.Item(2, 12).Value = [double]"0.02677421421684224"  # L2
$ws.Cells.Item(2, 13).Value = [double]"679.8868"  # M2
$ws.Cells.Item(2, 14).Value = [double]"0.04836126187004303"  # N2
$ws.Cells.Item(2, 15).Value = [double]"0.3968311260790642"  # O2
$ws.Cells.Item(2, 16).Value = [double]"547.8368"  # P2
$ws.Cells.Item(2, 17).Value = [double]"0.03896836789131131"  # Q2
$ws.Cells.Item(2, 18).Value = [double]"0.3197571923025291"  # R2
$ws.Cells.Item(2, 19).Value = [double]"132.05"  # S2
$ws.Cells.Item(2, 20).Value = [double]"0.194223508972376"  # T2
$ws.Cells.Item(2, 21).Value = [double]"0.091"  # U2
$ws.Cells.Item(2, 22).Value = [double]"6.47295230643383e-06"  # V2
$ws.Cells.Item(2, 23).Value = [double]"0.04933897365000786"  # W2
$ws.Cells.Item(2, 24).Value = [double]"0.05807901224726801"  # X2
$ws.Cells.Item(2, 25).Value = [double]"-0.00874003859726015"  # Y2
$ws.Cells.Item(2, 26).Value = [double]"2.14452423656267"  # Z2
$ws.Cells.Item(2, 27).Value = [double]"0.02888575347432119"  # AA2
$ws.Cells.Item(2, 28).Value = [double]"0.04334820328349445"  # AB2
$ws.Cells.Item(2, 29).Value = [double]"-0.0144268465946674"  # AC2
$ws.Cells.Item(2, 30).Value = [double]"4018.7"  # AD2
$ws.Cells.Item(2, 32).Value = [double]"4018.7"  # AF2
$ws.Cells.Item(2, 33).Value = [double]"4018.609"  # AG2
$ws.Cells.Item(2, 34).Value = [double]"0.2223076582656606"  # AH2
$ws.Cells.Item(2, 35).Value = [double]"0.1231989356092656"  # AI2
$ws.Cells.Item(2, 36).Value = [double]"0.2223037433695841"  # AJ2
$ws.Cells.Item(2, 37).Value = [double]"0.1231964895608944"  # AK2
$ws.Cells.Item(2, 38).Value = [double]"203.2"  # AL2
$ws.Cells.Item(2, 39).Value = [double]"203.2"  # AM2
$ws.Cells.Item(2, 40).Value = [double]"1.46979006656426"  # AN2
$ws.Cells.Item(2, 41).Value = [double]"10.42997047244094"  # AO2
$ws.Cells.Item(2, 42).Value = [double]"1.469756784434204"  # AP2
$ws.Cells.Item(2, 43).Value = [double]"10.42997047244094"  # AQ2
$ws.Cells.Item(3, 2).Value = "DB Insurance Co., Ltd. (KOSE:A005830)"  # B3
$ws.Cells.Item(3, 4).Value = [double]"0.0472"  # D3
$ws.Cells.Item(3, 5).Value = [double]"0.0457"  # E3
$ws.Cells.Item(3, 6).Value = [double]"0.0661"  # F3
$ws.Cells.Item(3, 7).Value = [double]"0.02567530880955613"  # G3
$ws.Cells.Item(3, 8).Value = [double]"0.02567530880955613"  # H3
$ws.Cells.Item(3, 9).Value = [double]"0.05264693905253157"  # I3
$ws.Cells.Item(3, 10).Value = [double]"0.03951002214250326"  # J3
$ws.Cells.Item(3, 11).Value = [double]"477.4"  # K3
$ws.Cells.Item(3, 12).Value = [double]"0.03240124881227094"  # L3
$ws.Cells.Item(3, 13).Value = [double]"182.9"  # M3
$ws.Cells.Item(3, 14).Value = [double]"0.07576322439004184"  # N3
$ws.Cells.Item(3, 15).Value = [double]"0.3831168831168831"  # O3
$ws.Cells.Item(3, 16).Value = [double]"77.4"  # P3
$ws.Cells.Item(3, 17).Value = [double]"0.03206163787746987"  # Q3
$ws.Cells.Item(3, 18).Value = [double]"0.1621281943862589"  # R3
$ws.Cells.Item(3, 19).Value = [double]"105.5"  # S3
$ws.Cells.Item(3, 20).Value = [double]"0.5768179332968836"  # T3
$ws.Cells.Item(3, 21).Value = [double]"0.0"  # U3
$ws.Cells.Item(3, 22).Value = [double]"0.0"  # V3
$ws.Cells.Item(3, 23).Value = [double]"0.08936058700209644"  # W3
$ws.Cells.Item(3, 24).Value = [double]"0.05807901224726801"  # X3
$ws.Cells.Item(3, 25).Value = [double]"0.03128157475482843"  # Y3
$ws.Cells.Item(3, 26).Value = [double]"2.339918688857833"  # Z3
$ws.Cells.Item(3, 27).Value = [double]"0.09245023920843015"  # AA3
$ws.Cells.Item(3, 28).Value = [double]"0.04330306607703813"  # AB3
$ws.Cells.Item(3, 29).Value = [double]"0.04914717313139202"  # AC3
$ws.Cells.Item(3, 30).Value = [double]"1392.6"  # AD3
$ws.Cells.Item(3, 32).Value = [double]"1392.6"  # AF3
$ws.Cells.Item(3, 33).Value = [double]"1392.6"  # AG3
$ws.Cells.Item(3, 34).Value = [double]"0.3658286705020096"  # AH3
$ws.Cells.Item(3, 35).Value = [double]"0.195310089478556"  # AI3
$ws.Cells.Item(3, 36).Value = [double]"0.3658286705020096"  # AJ3
$ws.Cells.Item(3, 37).Value = [double]"0.195310089478556"  # AK3
$ws.Cells.Item(3, 38).Value = [double]"44.0"  # AL3
$ws.Cells.Item(3, 39).Value = [double]"44.0"  # AM3
$ws.Cells.Item(3, 40).Value = [double]"1.584841242744964"  # AN3
$ws.Cells.Item(3, 41).Value = [double]"17.62954545454545"  # AO3
$ws.Cells.Item(3, 42).Value = [double]"1.584841242744964"  # AP3
$ws.Cells.Item(3, 43).Value = [double]"17.62954545454545"  # AQ3
$ws.Cells.Item(4, 2).Value = "Samsung Fire & Marine Insurance Co., Ltd. (KOSE:A000810)"  # B4
$ws.Cells.Item(4, 4).Value = [double]"0.00782"  # D4
$ws.Cells.Item(4, 5).Value = [double]"-0.0399"  # E4
$ws.Cells.Item(4, 6).Value = [double]"0.0576"  # F4
$ws.Cells.Item(4, 7).Value = [double]"0.05247635513068239"  # G4
$ws.Cells.Item(4, 8).Value = [double]"0.05247635513068239"  # H4
$ws.Cells.Item(4, 9).Value = [double]"0.0479695119013847"  # I4
$ws.Cells.Item(4, 10).Value = [double]"0.03541325684241432"  # J4
$ws.Cells.Item(4, 11).Value = [double]"596.0"  # K4
$ws.Cells.Item(4, 12).Value = [double]"0.03190117060167963"  # L4
$ws.Cells.Item(4, 13).Value = [double]"311.2"  # M4
$ws.Cells.Item(4, 14).Value = [double]"0.04248522164125108"  # N4
$ws.Cells.Item(4, 15).Value = [double]"0.5221476510067113"  # O4
$ws.Cells.Item(4, 16).Value = [double]"311.2"  # P4
$ws.Cells.Item(4, 17).Value = [double]"0.04248522164125108"  # Q4
$ws.Cells.Item(4, 18).Value = [double]"0.5221476510067113"  # R4
$ws.Cells.Item(4, 20).Value = [double]"0.0"  # T4
$ws.Cells.Item(4, 21).Value = [double]"0.086"  # U4
$ws.Cells.Item(4, 22).Value = [double]"1.174077461808352e-05"  # V4
$ws.Cells.Item(4, 23).Value = [double]"0.04933897365000786"  # W4
$ws.Cells.Item(4, 24).Value = [double]"0.04334820328349445"  # X4
$ws.Cells.Item(4, 25).Value = [double]"0.005990770366513409"  # Y4
$ws.Cells.Item(4, 26).Value = [double]"1.546625552902886"  # Z4
$ws.Cells.Item(4, 27).Value = [double]"0.05477104794399096"  # AA4
$ws.Cells.Item(4, 28).Value = [double]"0.04334820328349445"  # AB4
$ws.Cells.Item(4, 29).Value = [double]"0.01142284466049651"  # AC4
$ws.Cells.Item(4, 30).Value = [double]"0.0"  # AD4
$ws.Cells.Item(4, 32).Value = [double]"0.0"  # AF4
$ws.Cells.Item(4, 33).Value = [double]"-0.086"  # AG4
$ws.Cells.Item(4, 34).Value = [double]"0.0"  # AH4
$ws.Cells.Item(4, 35).Value = [double]"0.0"  # AI4
$ws.Cells.Item(4, 36).Value = [double]"-1.174091246549059e-05"  # AJ4
$ws.Cells.Item(4, 37).Value = [double]"-6.581508074515528e-06"  # AK4
$ws.Cells.Item(4, 38).Value = [double]"38.6"  # AL4
$ws.Cells.Item(4, 39).Value = [double]"38.6"  # AM4
$ws.Cells.Item(4, 40).Value = [double]"0.0"  # AN4
$ws.Cells.Item(4, 41).Value = [double]"23.21761658031088"  # AO4
$ws.Cells.Item(4, 42).Value = [double]"-7.321018132289095e-05"  # AP4
$ws.Cells.Item(4, 43).Value = [double]"23.21761658031088"  # AQ4
$ws.Cells.Item(5, 4).Value = [double]"0.0379"  # D5
$ws.Cells.Item(5, 5).Value = [double]"0.0356"  # E5
$ws.Cells.Item(5, 6).Value = [double]"0.006"  # F5
$ws.Cells.Item(5, 7).Value = [double]"0.02918817456593149"  # G5
$ws.Cells.Item(5, 8).Value = [double]"0.02918817456593149"  # H5
$ws.Cells.Item(5, 9).Value = [double]"0.02608321601751916"  # I5
$ws.Cells.Item(5, 10).Value = [double]"0.01883386419732224"  # J5
$ws.Cells.Item(5, 11).Value = [double]"312.9"  # K5
$ws.Cells.Item(5, 12).Value = [double]"0.02447207883622712"  # L5
$ws.Cells.Item(5, 13).Value = [double]"81.4264"  # M5
$ws.Cells.Item(5, 14).Value = [double]"0.0490431849665723"  # N5
$ws.Cells.Item(5, 15).Value = [double]"0.2602313838286993"  # O5
$ws.Cells.Item(5, 16).Value = [double]"60.0264"  # P5
$ws.Cells.Item(5, 17).Value = [double]"0.03615394808167199"  # Q5
$ws.Cells.Item(5, 18).Value = [double]"0.1918389261744967"  # R5
$ws.Cells.Item(5, 19).Value = [double]"21.4"  # S5
$ws.Cells.Item(5, 20).Value = [double]"0.2628140259178841"  # T5
$ws.Cells.Item(5, 23).Value = [double]"0.07772368224949078"  # W5
$ws.Cells.Item(5, 24).Value = [double]"0.05712139158095693"  # X5
$ws.Cells.Item(5, 25).Value = [double]"0.02060229066853385"  # Y5
$ws.Cells.Item(5, 26).Value = [double]"2.574759862260617"  # Z5
$ws.Cells.Item(5, 27).Value = [double]"0.04849267758653257"  # AA5
$ws.Cells.Item(5, 28).Value = [double]"0.04330497224555261"  # AB5
$ws.Cells.Item(5, 29).Value = [double]"0.00518770534097996"  # AC5
$ws.Cells.Item(5, 30).Value = [double]"895.5"  # AD5
$ws.Cells.Item(5, 32).Value = [double]"895.5"  # AF5
$ws.Cells.Item(5, 33).Value = [double]"895.5"  # AG5
$ws.Cells.Item(5, 34).Value = [double]"0.3503795289146255"  # AH5
$ws.Cells.Item(5, 35).Value = [double]"0.172099012184341"  # AI5
$ws.Cells.Item(5, 36).Value = [double]"0.3503795289146255"  # AJ5
$ws.Cells.Item(5, 37).Value = [double]"0.172099012184341"  # AK5
$ws.Cells.Item(5, 38).Value = [double]"36.6"  # AL5
$ws.Cells.Item(5, 39).Value = [double]"36.6"  # AM5
$ws.Cells.Item(5, 40).Value = [double]"2.106563161609033"  # AN5
$ws.Cells.Item(5, 41).Value = [double]"9.112021857923496"  # AO5
$ws.Cells.Item(5, 42).Value = [double]"2.106563161609033"  # AP5
$ws.Cells.Item(5, 43).Value = [double]"9.112021857923496"  # AQ5
$ws.Cells.Item(6, 2).Value = "Meritz Fire & Marine Insurance Co., Ltd. (KOSE:A000060)"  # B6
$ws.Cells.Item(6, 4).Value = [double]"0.09720000000000001"  # D6
$ws.Cells.Item(6, 5).Value = [double]"0.22"  # E6
$ws.Cells.Item(6, 6).Value = [double]"0.006"  # F6
$ws.Cells.Item(6, 7).Value = [double]"-0.02831614782405644"  # G6
$ws.Cells.Item(6, 8).Value = [double]"-0.02831614782405644"  # H6
$ws.Cells.Item(6, 9).Value = [double]"0.01453906374820259"  # I6
$ws.Cells.Item(6, 10).Value = [double]"0.01045171385776649"  # J6
$ws.Cells.Item(6, 11).Value = [double]"332.8"  # K6
$ws.Cells.Item(6, 12).Value = [double]"0.04090110241252597"  # L6
$ws.Cells.Item(6, 13).Value = [double]"91.29"  # M6
$ws.Cells.Item(6, 14).Value = [double]"0.05764713311442284"  # N6
$ws.Cells.Item(6, 15).Value = [double]"0.2743088942307693"  # O6
$ws.Cells.Item(6, 16).Value = [double]"86.14"  # P6
$ws.Cells.Item(6, 17).Value = [double]"0.05439504925486234"  # Q6
$ws.Cells.Item(6, 18).Value = [double]"0.2588341346153846"  # R6
$ws.Cells.Item(6, 19).Value = [double]"5.150000000000006"  # S6
$ws.Cells.Item(6, 20).Value = [double]"0.05641362690327534"  # T6
$ws.Cells.Item(6, 21).Value = [double]"0.003"  # U6
$ws.Cells.Item(6, 22).Value = [double]"1.89441778226825e-06"  # V6
$ws.Cells.Item(6, 23).Value = [double]"0.1378225038307036"  # W6
$ws.Cells.Item(6, 24).Value = [double]"0.05370553513922909"  # X6
$ws.Cells.Item(6, 25).Value = [double]"0.08411696869147454"  # Y6
$ws.Cells.Item(6, 26).Value = [double]"2.763733667742605"  # Z6
$ws.Cells.Item(6, 27).Value = [double]"0.02888575347432119"  # AA6
$ws.Cells.Item(6, 28).Value = [double]"0.04331260006898859"  # AB6
$ws.Cells.Item(6, 29).Value = [double]"-0.0144268465946674"  # AC6
$ws.Cells.Item(6, 30).Value = [double]"642.3"  # AD6
$ws.Cells.Item(6, 32).Value = [double]"642.3"  # AF6
$ws.Cells.Item(6, 33).Value = [double]"642.2969999999999"  # AG6
$ws.Cells.Item(6, 34).Value = [double]"0.2885574374410351"  # AH6
$ws.Cells.Item(6, 35).Value = [double]"0.2146796350145393"  # AI6
$ws.Cells.Item(6, 36).Value = [double]"0.2885564785791975"  # AJ6
$ws.Cells.Item(6, 37).Value = [double]"0.2146788475672792"  # AK6
$ws.Cells.Item(6, 38).Value = [double]"33.3"  # AL6
$ws.Cells.Item(6, 39).Value = [double]"33.3"  # AM6
$ws.Cells.Item(6, 40).Value = [double]"3.873944511459589"  # AN6
$ws.Cells.Item(6, 41).Value = [double]"3.552552552552553"  # AO6
$ws.Cells.Item(6, 42).Value = [double]"3.873926417370325"  # AP6
$ws.Cells.Item(6, 43).Value = [double]"3.552552552552553"  # AQ6
$ws.Cells.Item(7, 2).Value = "Hanwha General Insurance Co., Ltd. (KOSE:A000370)"  # B7
$ws.Cells.Item(7, 4).Value = [double]"0.0199"  # D7
$ws.Cells.Item(7, 5).Value = [double]"-0.269"  # E7
$ws.Cells.Item(7, 7).Value = [double]"-0.006668519033064741"  # G7
$ws.Cells.Item(7, 8).Value = [double]"-0.006668519033064741"  # H7
$ws.Cells.Item(7, 9).Value = [double]"0.005749460320174407"  # I7
$ws.Cells.Item(7, 10).Value = [double]"0.005749460320174407"  # J7
$ws.Cells.Item(7, 11).Value = [double]"-3.61"  # K7
$ws.Cells.Item(7, 12).Value = [double]"-0.0007715818496590933"  # L7
$ws.Cells.Item(7, 13).Value = [double]"13.0704"  # M7
$ws.Cells.Item(7, 14).Value = [double]"0.031770539620807"  # N7
$ws.Cells.Item(7, 15).Value = [double]"-3.620609418282549"  # O7
$ws.Cells.Item(7, 16).Value = [double]"13.0704"  # P7
$ws.Cells.Item(7, 17).Value = [double]"0.031770539620807"  # Q7
$ws.Cells.Item(7, 18).Value = [double]"-3.620609418282549"  # R7
$ws.Cells.Item(7, 20).Value = [double]"0.0"  # T7
$ws.Cells.Item(7, 21).Value = [double]"0.0"  # U7
$ws.Cells.Item(7, 22).Value = [double]"0.0"  # V7
$ws.Cells.Item(7, 23).Value = [double]"-0.002791309054357071"  # W7
$ws.Cells.Item(7, 24).Value = [double]"0.0689091842852136"  # X7
$ws.Cells.Item(7, 25).Value = [double]"-0.07170049333957067"  # Y7
$ws.Cells.Item(7, 26).Value = [double]"2.758667452830188"  # Z7
$ws.Cells.Item(7, 27).Value = [double]"0.01586084905660377"  # AA7
$ws.Cells.Item(7, 28).Value = [double]"0.043997424052615"  # AB7
$ws.Cells.Item(7, 29).Value = [double]"-0.02813657499601122"  # AC7
$ws.Cells.Item(7, 30).Value = [double]"411.8"  # AD7
$ws.Cells.Item(7, 32).Value = [double]"411.8"  # AF7
$ws.Cells.Item(7, 33).Value = [double]"411.8"  # AG7
$ws.Cells.Item(7, 34).Value = [double]"0.500242954324587"  # AH7
$ws.Cells.Item(7, 35).Value = [double]"0.198841139546113"  # AI7
$ws.Cells.Item(7, 36).Value = [double]"0.500242954324587"  # AJ7
$ws.Cells.Item(7, 37).Value = [double]"0.198841139546113"  # AK7
$ws.Cells.Item(7, 38).Value = [double]"19.0"  # AL7
$ws.Cells.Item(7, 39).Value = [double]"19.0"  # AM7
$ws.Cells.Item(7, 40).Value = [double]"5.040391676866585"  # AN7
$ws.Cells.Item(7, 41).Value = [double]"1.41578947368421"  # AO7
$ws.Cells.Item(7, 42).Value = [double]"5.040391676866585"  # AP7
$ws.Cells.Item(7, 43).Value = [double]"1.41578947368421"  # AQ7
$ws.Cells.Item(8, 2).Value = "Heungkuk Fire & Marine Insurance Co., Ltd. (KOSE:A000540)"  # B8
$ws.Cells.Item(8, 4).Value = [double]"-0.0111"  # D8
$ws.Cells.Item(8, 5).Value = [double]"0.102"  # E8
$ws.Cells.Item(8, 7).Value = [double]"0.01904030054644809"  # G8
$ws.Cells.Item(8, 8).Value = [double]"0.01904030054644809"  # H8
$ws.Cells.Item(8, 9).Value = [double]"0.002634050546448087"  # I8
$ws.Cells.Item(8, 10).Value = [double]"0.002634050546448087"  # J8
$ws.Cells.Item(8, 11).Value = [double]"20.9"  # K8
$ws.Cells.Item(8, 12).Value = [double]"0.008922472677595628"  # L8
$ws.Cells.Item(8, 13).Value = [double]"-0.0"  # M8
$ws.Cells.Item(8, 14).Value = [double]"-0.0"  # N8
$ws.Cells.Item(8, 15).Value = [double]"-0.0"  # O8
$ws.Cells.Item(8, 16).Value = [double]"-0.0"  # P8
$ws.Cells.Item(8, 17).Value = [double]"-0.0"  # Q8
$ws.Cells.Item(8, 18).Value = [double]"-0.0"  # R8
$ws.Cells.Item(8, 21).Value = [double]"0.002"  # U8
$ws.Cells.Item(8, 22).Value = [double]"1.196888090963495e-05"  # V8
$ws.Cells.Item(8, 23).Value = [double]"0.03334396936821953"  # W8
$ws.Cells.Item(8, 24).Value = [double]"0.09235744382665725"  # X8
$ws.Cells.Item(8, 25).Value = [double]"-0.05901347445843773"  # Y8
$ws.Cells.Item(8, 26).Value = [double]"2.368270992320105"  # Z8
$ws.Cells.Item(8, 27).Value = [double]"0.006238145501457927"  # AA8
$ws.Cells.Item(8, 28).Value = [double]"0.04652142166187055"  # AB8
$ws.Cells.Item(8, 29).Value = [double]"-0.04028327616041263"  # AC8
$ws.Cells.Item(8, 30).Value = [double]"320.7"  # AD8
$ws.Cells.Item(8, 32).Value = [double]"320.7"  # AF8
$ws.Cells.Item(8, 33).Value = [double]"320.698"  # AG8
$ws.Cells.Item(8, 34).Value = [double]"0.6574415744157442"  # AH8
$ws.Cells.Item(8, 35).Value = [double]"0.3413881200766446"  # AI8
$ws.Cells.Item(8, 36).Value = [double]"0.6574401699063956"  # AJ8
$ws.Cells.Item(8, 37).Value = [double]"0.3413867178767678"  # AK8
$ws.Cells.Item(8, 38).Value = [double]"17.6"  # AL8
$ws.Cells.Item(8, 39).Value = [double]"17.6"  # AM8
$ws.Cells.Item(8, 40).Value = [double]"12.67588932806324"  # AN8
$ws.Cells.Item(8, 41).Value = [double]"0.3505681818181818"  # AO8
$ws.Cells.Item(8, 42).Value = [double]"12.67581027667984"  # AP8
$ws.Cells.Item(8, 43).Value = [double]"0.3505681818181818"  # AQ8
$ws.Cells.Item(9, 2).Value = "Lotte Non - Life Insurance Co., Ltd. (KOSE:A000400)"  # B9
$ws.Cells.Item(9, 4).Value = [double]"0.0897"  # D9
$ws.Cells.Item(9, 7).Value = [double]"-0.02680812229066849"  # G9
$ws.Cells.Item(9, 8).Value = [double]"-0.02680812229066849"  # H9
$ws.Cells.Item(9, 9).Value = [double]"-0.01422161381093619"  # I9
$ws.Cells.Item(9, 10).Value = [double]"-0.01422161381093619"  # J9
$ws.Cells.Item(9, 11).Value = [double]"-23.1"  # K9
$ws.Cells.Item(9, 12).Value = [double]"-0.008783937942048825"  # L9
$ws.Cells.Item(9, 13).Value = [double]"-0.0"  # M9
$ws.Cells.Item(9, 14).Value = [double]"-0.0"  # N9
$ws.Cells.Item(9, 15).Value = [double]"0.0"  # O9
$ws.Cells.Item(9, 16).Value = [double]"-0.0"  # P9
$ws.Cells.Item(9, 17).Value = [double]"-0.0"  # Q9
$ws.Cells.Item(9, 18).Value = [double]"0.0"  # R9
$ws.Cells.Item(9, 19).Value = [double]"0.0"  # S9
$ws.Cells.Item(9, 21).Value = [double]"0.0"  # U9
$ws.Cells.Item(9, 22).Value = [double]"0.0"  # V9
$ws.Cells.Item(9, 23).Value = [double]"-0.03844233649525712"  # W9
$ws.Cells.Item(9, 24).Value = [double]"0.06162573908353568"  # X9
$ws.Cells.Item(9, 25).Value = [double]"-0.1000680755787928"  # Y9
$ws.Cells.Item(9, 26).Value = [double]"3.031819229882407"  # Z9
$ws.Cells.Item(9, 27).Value = [double]"-0.04311736223195757"  # AA9
$ws.Cells.Item(9, 28).Value = [double]"0.04536169861189269"  # AB9
$ws.Cells.Item(9, 29).Value = [double]"-0.08847906084385027"  # AC9
$ws.Cells.Item(9, 30).Value = [double]"355.8"  # AD9
$ws.Cells.Item(9, 32).Value = [double]"355.8"  # AF9
$ws.Cells.Item(9, 33).Value = [double]"355.8"  # AG9
$ws.Cells.Item(9, 34).Value = [double]"0.4171649665845937"  # AH9
$ws.Cells.Item(9, 35).Value = [double]"0.2924303443741267"  # AI9
$ws.Cells.Item(9, 36).Value = [double]"0.4171649665845937"  # AJ9
$ws.Cells.Item(9, 37).Value = [double]"0.2924303443741267"  # AK9
$ws.Cells.Item(9, 38).Value = [double]"14.1"  # AL9
$ws.Cells.Item(9, 39).Value = [double]"14.1"  # AM9
$ws.Cells.Item(9, 40).Value = [double]"-20.80701754385965"  # AN9
$ws.Cells.Item(9, 41).Value = [double]"-2.652482269503546"  # AO9
$ws.Cells.Item(9, 42).Value = [double]"-20.80701754385965"  # AP9
$ws.Cells.Item(9, 43).Value = [double]"-2.652482269503546"  # AQ9

$ws.Cells.Item(8, 6).ClearContents()  # F8
$ws.Cells.Item(8, 20).ClearContents()  # T8
$ws.Cells.Item(9, 6).ClearContents()  # F9
$ws.Cells.Item(9, 20).ClearContents()  # T9
